$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (and a row reorder for rows 36-37)
# Each cell is written with a leading quote-prefix so Excel stores the
# value as literal text (matching the source inlineStr cells) instead of
# auto-converting numeric-looking strings into numbers; the style is then
# reset to "Normal" so no stray number-format style is left on the cell.

$ws.Range("D2").Value = "'96.431.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.92%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.705.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.02%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'238.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.38%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +8.99%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'654.73"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.26%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.421"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.31%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +1.76%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.04%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'3.702.84"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +3.05%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'44.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.82%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  +0.31%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +5.33%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.397.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.93%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.0000267"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.20%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'96.406.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.96%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'8.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +15.67%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.710.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.17%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'19.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +4.25%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'12.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.76%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.525"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.04%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'521.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.64%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'3.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.56%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'7.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.24%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0000203"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.07%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'101.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.42%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'13.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.30%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -5.91%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'12.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +4.24%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'3.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.10%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.08%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +11.94%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -3.00%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'32.67"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.70%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'Binance-PegBSC-USD"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.08%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'Bittensor"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'661.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +7.27%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.594"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.58%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'8.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.74%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'7.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +16.40%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'41.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +24.63%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.161"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +3.49%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.970"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +4.15%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.06%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +0.02%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.447"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.88%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.0455"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +1.97%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.94%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'23.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.19%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'8.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.42%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +1.77%  "
$ws.Range("E51").Style = "Normal"
